# Apply the "Add files via upload" revision to the SEO audit table:
# - refresh the audit rows (columns A-D) with updated wording/ordering
# - add a "Référence" hyperlink in column F for each row that gained one
# - row 20 becomes a blank spacer row (only keeps its new reference link)
# - move the active selection to E26, mirroring the saved workbook state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh audit rows (columns A-D) ---
# Row 3
$ws.Range("A3").Value = 'SEO'
$ws.Range("B3").Value = 'Keywords'
$ws.Range("C3").Value = 'Inutile dans le code'
$ws.Range("D3").Value = 'Les supprimer du HTML'

# Row 4
$ws.Range("A4").Value = 'Accessibilité'
$ws.Range("B4").Value = 'Couleurs illisibles'
$ws.Range("C4").Value = 'On ne peut pas lire'
$ws.Range("D4").Value = 'Changement de couleur'

# Row 5
$ws.Range("A5").Value = 'SEO'
$ws.Range("B5").Value = 'lang=default'
$ws.Range("C5").Value = 'Une langue prédéfini'
$ws.Range("D5").Value = 'Mettre sur "fr"'

# Row 6
$ws.Range("A6").Value = 'SEO '
$ws.Range("B6").Value = '.min'
$ws.Range("C6").Value = 'Fichier css non racordés'
$ws.Range("D6").Value = 'Enlever les .min'

# Row 7
$ws.Range("A7").Value = 'SEO'
$ws.Range("B7").Value = 'Pas de class aux labels'
$ws.Range("C7").Value = 'symbole rouge avec Wave'
$ws.Range("D7").Value = 'Ajouter des class aux label'

# Row 8
$ws.Range("A8").Value = 'SEO '
$ws.Range("B8").Value = 'Title'
$ws.Range("C8").Value = 'Pas de titre'
$ws.Range("D8").Value = 'Mettre un titre'

# Row 9
$ws.Range("A9").Value = 'SEO '
$ws.Range("B9").Value = 'Image non indispensable'
$ws.Range("C9").Value = 'Prise de poinds'
$ws.Range("D9").Value = 'Mettre un titre à la place'

# Row 10
$ws.Range("A10").Value = 'Accessibilité'
$ws.Range("B10").Value = 'Liens dans le footer'
$ws.Range("C10").Value = 'Aucune utilité'
$ws.Range("D10").Value = 'Les enlever'

# Row 11
$ws.Range("A11").Value = 'SEO '
$ws.Range("B11").Value = 'Entité spécifique en HTML'
$ws.Range("C11").Value = 'Aucune utilité'
$ws.Range("D11").Value = 'Les enlever et les remplacer'

# Row 12
$ws.Range("A12").Value = 'Accessibilité'
$ws.Range("B12").Value = 'Tailles des paragraphes'
$ws.Range("C12").Value = 'On ne peut pas lire'
$ws.Range("D12").Value = 'Augmenter la font-size'

# Row 13
$ws.Range("A13").Value = 'SEO'
$ws.Range("B13").Value = 'Les photos (images)'
$ws.Range("C13").Value = 'Trop lourdes'
$ws.Range("D13").Value = 'Baisser leurs poinds'

# Row 14
$ws.Range("A14").Value = 'SEO '
$ws.Range("B14").Value = 'Bootstrap et Jquery'
$ws.Range("C14").Value = 'Anciennes versions'
$ws.Range("D14").Value = 'Les remplacer'

# Row 15
$ws.Range("A15").Value = 'SEO '
$ws.Range("B15").Value = 'Les images'
$ws.Range("C15").Value = 'Le codec'
$ws.Range("D15").Value = 'Changer pour PNG ou JPG'

# Row 16
$ws.Range("A16").Value = 'SEO'
$ws.Range("B16").Value = 'les balises "script" (JS)'
$ws.Range("C16").Value = 'Placer au mauvais endroit'
$ws.Range("D16").Value = 'Les mettre à la fin du code (body)'

# Row 17
$ws.Range("A17").Value = 'Accessibilité'
$ws.Range("B17").Value = 'Les paragraphes'
$ws.Range("C17").Value = 'Trop rapprochés'
$ws.Range("D17").Value = 'Faire des espaces'

# Row 18
$ws.Range("A18").Value = 'SEO '
$ws.Range("B18").Value = 'Page2'
$ws.Range("C18").Value = 'Nom inexact'
$ws.Range("D18").Value = 'Changer par "Contact"'

# Row 19
$ws.Range("A19").Value = 'SEO'
$ws.Range("B19").Value = 'Width et height'
$ws.Range("C19").Value = 'Rien à faire dans le HTML'
$ws.Range("D19").Value = 'Remettre dans le CSS'

# Row 20
$ws.Range("A20:E20").ClearContents()

# Row 21
$ws.Range("A21").Value = 'SEO'
$ws.Range("B21").Value = '<li>….</li>'
$ws.Range("C21").Value = 'Il y en a en trop'
$ws.Range("D21").Value = 'Les supprimer'

# --- Add "Référence" hyperlinks in column F ---
# (F4 and F7 already show the plain "Wave Evaluation Tool" label and are left as-is)
$ws.Range("F3").Value = 'https://optimiz.me/la-balise-meta-keywords/'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://optimiz.me/la-balise-meta-keywords/') | Out-Null

$ws.Range("F5").Value = 'http://www.oujood.com/html-attribut/HTML_attribut-lang.php'
$ws.Hyperlinks.Add($ws.Range("F5"), 'http://www.oujood.com/html-attribut/HTML_attribut-lang.php') | Out-Null

$ws.Range("F6").Value = 'https://developer.mozilla.org/fr/docs/Apprendre/Commencer_avec_le_web/G%C3%A9rer_les_fichiers'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://developer.mozilla.org/fr/docs/Apprendre/Commencer_avec_le_web/G%C3%A9rer_les_fichiers') | Out-Null

$ws.Range("F8").Value = 'https://developer.mozilla.org/fr/docs/Web/HTML/Element/title'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://developer.mozilla.org/fr/docs/Web/HTML/Element/title') | Out-Null

$ws.Range("F9").Value = 'https://developer.mozilla.org/fr/docs/Web/HTML/Element/Heading_Elements'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://developer.mozilla.org/fr/docs/Web/HTML/Element/Heading_Elements') | Out-Null

$ws.Range("F10").Value = 'https://www.209-agency.com/actualites-seo/actus-seo/footer-et-seo-utile-pour-google/'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.209-agency.com/actualites-seo/actus-seo/footer-et-seo-utile-pour-google/') | Out-Null

$ws.Range("F11").Value = 'https://developer.mozilla.org/fr/docs/Glossaire/Entity'
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://developer.mozilla.org/fr/docs/Glossaire/Entity') | Out-Null

$ws.Range("F12").Value = 'https://developer.mozilla.org/fr/docs/Web/CSS/font-size'
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://developer.mozilla.org/fr/docs/Web/CSS/font-size') | Out-Null

$ws.Range("F13").Value = 'https://www.oscar-referencement.com/astuces-optimiser-images-referencement-naturel/'
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.oscar-referencement.com/astuces-optimiser-images-referencement-naturel/') | Out-Null

$ws.Range("F14").Value = 'https://www.developpez.com/actu/246027/La-version-5-0-du-framework-Bootstrap-va-supprimer-jQuery-sa-plus-grande-dependance-cote-client-pour-du-pur-JavaScript/'
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.developpez.com/actu/246027/La-version-5-0-du-framework-Bootstrap-va-supprimer-jQuery-sa-plus-grande-dependance-cote-client-pour-du-pur-JavaScript/') | Out-Null

$ws.Range("F15").Value = 'https://developer.mozilla.org/fr/docs/Web/Media/Formats/Types_des_images'
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://developer.mozilla.org/fr/docs/Web/Media/Formats/Types_des_images') | Out-Null

$ws.Range("F16").Value = 'https://initiativedeveloppeur.fr/pourquoi-placer-javascript-fin-et-css-debut/'
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://initiativedeveloppeur.fr/pourquoi-placer-javascript-fin-et-css-debut/') | Out-Null

$ws.Range("F17").Value = 'https://developer.mozilla.org/fr/docs/Web/CSS/margin'
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://developer.mozilla.org/fr/docs/Web/CSS/margin') | Out-Null

$ws.Range("F18").Value = 'https://developer.mozilla.org/fr/docs/Web/HTML/Element/a'
$ws.Hyperlinks.Add($ws.Range("F18"), 'https://developer.mozilla.org/fr/docs/Web/HTML/Element/a') | Out-Null

$ws.Range("F19").Value = 'https://developer.mozilla.org/fr/docs/Web/CSS/height'
$ws.Hyperlinks.Add($ws.Range("F19"), 'https://developer.mozilla.org/fr/docs/Web/CSS/height') | Out-Null

$ws.Range("F20").Value = 'https://developer.mozilla.org/fr/docs/Web/CSS/width'
$ws.Hyperlinks.Add($ws.Range("F20"), 'https://developer.mozilla.org/fr/docs/Web/CSS/width') | Out-Null

$ws.Range("F21").Value = 'https://developer.mozilla.org/fr/docs/Web/HTML/Element/li'
$ws.Hyperlinks.Add($ws.Range("F21"), 'https://developer.mozilla.org/fr/docs/Web/HTML/Element/li') | Out-Null

# --- Restore the saved selection/active cell ---
$ws.Range("E26").Select()
